# Ajout synthese des chiffres sur recuit
#
# The sheet holds 27 blocks of results. Each block is a label row (column A,
# shared-string) followed by 10 data rows (columns A and B). This change adds,
# on the first data row of every block, two summary formulas in columns D/E:
#   D<row> = AVERAGE(A<row>:A<row+9>)
#   E<row> = AVERAGE(B<row>:B<row+9>)
# (Excel's fill-right relative-reference adjustment turns the single formula
# typed into D<row>:E<row> into the A-column average in D and the B-column
# average in E.)
#
# Two data fixes are also part of this edit:
#   - row 276 (a spurious 11th data row in the "-- 50-10/ --" block) is
#     cleared out entirely so that block also has exactly 10 data rows;
#   - a new row 298 (A298=0, B298=0) is added so the last block
#     ("-- 50-1000/ --") also has exactly 10 data rows (it only had 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data fix #1: drop the spurious extra row in the "50-10" block ---
$ws.Range("A276").ClearContents()
$ws.Range("B276").ClearContents()

# --- data fix #2: complete the last block ("50-1000") with a 10th row ---
$ws.Range("A298").Value = 0
$ws.Range("B298").Value = 0

# --- add the AVERAGE summary formulas for every block ---
$starts = @(2, 13, 24, 35, 46, 57, 68, 79, 90, 101, 112, 123, 134, 145, 156, 167, 178, 189, 200, 211, 222, 233, 244, 255, 266, 278, 289)
foreach ($start in $starts) {
    $end = $start + 9
    $ws.Range("D${start}:E${start}").Formula = "=AVERAGE(A${start}:A${end})"
}

# --- restore the view state (scroll position + selection) ---
$excel.ActiveWindow.ScrollRow = 265
$ws.Range("D289").Select()
